$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted at row 849 on the "Feria Lagunitas de
# Puerto Montt - Lechuga" sheet. Inserting a whole row shifts every
# subsequent row (old 849..936) down by one (new 850..937), which matches
# the commit's weekly-update pattern.
$ws.Rows.Item(849).Insert()

$ws.Cells.Item(849, 1).Value  = 4
$ws.Cells.Item(849, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(849, 3).Value  = "Los Lagos"
$ws.Cells.Item(849, 4).Value  = 45132
$ws.Cells.Item(849, 5).Value  = 10
$ws.Cells.Item(849, 6).Value  = 100112033
$ws.Cells.Item(849, 7).Value  = "Lechuga"
$ws.Cells.Item(849, 8).Value  = "Escarola"
$ws.Cells.Item(849, 9).Value  = "Primera"
$ws.Cells.Item(849, 10).Value = 400
$ws.Cells.Item(849, 11).Value = 13000
$ws.Cells.Item(849, 12).Value = 13000
$ws.Cells.Item(849, 13).Value = 13000
$ws.Cells.Item(849, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(849, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(849, 16).Value = 867
$ws.Cells.Item(849, 17).Value = 15
$ws.Cells.Item(849, 18).Value = "Hortaliza"
